$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: A="A" (shared string idx 4), B=date, C="B" (shared string idx 5), D=date
$ws.Range("A3").Value = "A"
$ws.Range("B3").Value = 43970.748611111114
$ws.Range("C3").Value = "B"
$ws.Range("D3").Value = 43971.848611111112

# Row 4: A="C" (new shared string idx 6), B=date, C="A" (shared string idx 4), D=date
$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = 43971.815972222219
$ws.Range("C4").Value = "A"
$ws.Range("D4").Value = 43972.640277777777

# Apply the same date/time number format as the existing date cells (B2/D2) to the new cells
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D2").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F6").Select()
